# "updated test cases with global share value for login user"
#
# Rows 45-48 (columns B/F/G/H) of the SystemTest execution-result sheet get
# new status values. We reuse the formatting of existing "donor" cells that
# already carry the right cell style (bold-red "Failed", green "Passed",
# gold "OnHold"/"Report Hidden", bold-red "M & S Needed") so the workbook's
# style table (cellXfs) is reused rather than duplicated, then overwrite the
# pasted value with the correct text for each target cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

function Copy-StyleAndSet($donorAddr, $targetAddr, $value) {
    $ws.Range($donorAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial($xlPasteFormats) | Out-Null
    if ($null -ne $value) {
        $ws.Range($targetAddr).Value2 = $value
    }
}

# Row 45 -------------------------------------------------------------
# B45: style only changes (still "Failed") -> pick up the red/bold style
# used elsewhere for "Failed" in column B/F (e.g. B17).
Copy-StyleAndSet "B17" "B45" "Failed"

# F45: was blank -> "Failed" with the same red/bold style.
Copy-StyleAndSet "F17" "F45" "Failed"

# G45: was blank -> new free-text note "Sudip will send" (plain style).
$ws.Range("G45").Value2 = "Sudip will send"

# Row 46 -------------------------------------------------------------
# F46: was blank -> "Passed" with the green style used in column B/F (e.g. F8).
Copy-StyleAndSet "F8" "F46" "Passed"

# Row 47 -------------------------------------------------------------
# F47/G47/H47: was blank -> "OnHold" / (blank) / "Report Hidden" with the
# gold style used elsewhere for that trio (e.g. row 23).
Copy-StyleAndSet "F23" "F47" "OnHold"
Copy-StyleAndSet "G23" "G47" $null
Copy-StyleAndSet "H23" "H47" "Report Hidden"

# Row 48 -------------------------------------------------------------
# F48: was blank -> "OnHold" (gold style, e.g. row 34).
Copy-StyleAndSet "F34" "F48" "OnHold"
# H48: was blank -> "M & S Needed" (red/bold style, e.g. H34).
Copy-StyleAndSet "H34" "H48" "M & S Needed"

# View state -----------------------------------------------------------
# Move the active selection to K42 (scrolled so row 28 is at the top).
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K42").Select() | Out-Null
